$d = $word.ActiveDocument

# The "Things for the immediate future" checklist block currently reads
# (paragraph indices 24-29):
#   24 Navbar                                                            (ilvl 2)
#   25 Make the search bar appear more uniform  [has the _GoBack bookmark] (ilvl 3)
#   26 Sidebar:                                                          (ilvl 2)
#   27 Scaling                                                           (ilvl 3)
#   28 Set to active when screen goes below a threshold, if the user...  (ilvl 4)
#   29 Collapse sidebar and center content based on sidebar width        (ilvl 3)
#
# It needs to become:
#   Sidebar:                                                             (ilvl 2)
#   Collapse sidebar and center content based on sidebar width           (ilvl 3)
#   Bug when you click on the button and an overlay comes on top and
#     then you make the window larger, the overlay stays and when you
#     click the button, it affects the effect of having the sidebar
#     disappear                                                         (ilvl 3)
#   Navbar                                                               (ilvl 2)
#   On mobile sized screens, the navbar collapses                       (ilvl 3)
#
# Work from the bottom of the block upward so deletions don't disturb the
# indices of paragraphs we haven't processed yet.

# 1. Drop the old "Collapse sidebar..." paragraph (#29) and the
#    "Set to active..." paragraph (#28) - their content is replaced by
#    fresh paragraphs inserted later in the block.
$p28 = $d.Paragraphs.Item(28)
$p29 = $d.Paragraphs.Item(29)
$deadRange = $d.Range($p28.Range.Start, $p29.Range.End)
$deadRange.Delete()

# 2. "Scaling" (#27) becomes "On mobile sized screens, the navbar collapses".
$d.Paragraphs.Item(27).Range.Text = "On mobile sized screens, the navbar collapses"

# 3. Drop the old "Sidebar:" paragraph (#26) - it is recreated at the top
#    of the block from paragraph #24 below.
$d.Paragraphs.Item(26).Range.Delete()

# 4. "Make the search bar appear more uniform" (#25, carries the
#    _GoBack bookmark) becomes "Navbar" and moves up one list level.
$p25 = $d.Paragraphs.Item(25)
$p25.Range.Text = "Navbar"
$p25.Range.ListFormat.ListLevelNumber = 3

# 5. "Navbar" (#24) becomes "Sidebar:".
$p24 = $d.Paragraphs.Item(24)
$p24.Range.Text = "Sidebar:"

# 6. Insert the two new sidebar bullets right after the renamed
#    "Sidebar:" paragraph, at ilvl 3 (ListLevelNumber 4), matching the
#    level used by the other bullets under "Sidebar:".
$p24.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs.Item(25)
$newPara1.Range.ListFormat.ListLevelNumber = 4
$newPara1.Range.Text = "Collapse sidebar and center content based on sidebar width"

$newPara1.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item(26)
$newPara2.Range.ListFormat.ListLevelNumber = 4
$newPara2.Range.Text = "Bug when you click on the button and an overlay comes on top and then you make the window larger, the overlay stays and when you click the button, it affects the effect of having the sidebar disappear"
